# Add five new custom styles used by the "elegant" CV template:
#   ContactInfo, SkillCategory, SkillItems (paragraph styles)
#   SkillHighlight, SkillLevel (character styles)
#
# Word's Font.Color takes a decimal value interpreted as 0xBBGGRR, so
# convert our target RGB hex colors accordingly.
$colorBrown = [Convert]::ToInt32("122D7C", 16)   # RGB 7C2D12 (body/contact text)
$colorAmber = [Convert]::ToInt32("0677D9", 16)   # RGB D97706 (highlight/category)

$d = $word.ActiveDocument

# --- ContactInfo (paragraph) ---
$contactInfo = $d.Styles.Add("ContactInfo", 1)
$contactInfo.NameLocal = "Contact Info"
$contactInfo.ParagraphFormat.SpaceBefore = 0
$contactInfo.ParagraphFormat.SpaceAfter = 5.1
$contactInfo.ParagraphFormat.LineSpacingRule = 0
$contactInfo.ParagraphFormat.LineSpacing = 18
$contactInfo.Font.Name = "Liberation Serif"
$contactInfo.Font.Bold = $false
$contactInfo.Font.Color = $colorBrown
$contactInfo.Font.Size = 9

# --- SkillCategory (paragraph) ---
$skillCategory = $d.Styles.Add("SkillCategory", 1)
$skillCategory.NameLocal = "Skill Category"
$skillCategory.ParagraphFormat.SpaceBefore = 0
$skillCategory.ParagraphFormat.SpaceAfter = 0
$skillCategory.ParagraphFormat.LineSpacingRule = 0
$skillCategory.ParagraphFormat.LineSpacing = 14.4
$skillCategory.Font.Name = "Liberation Serif"
$skillCategory.Font.Bold = $true
$skillCategory.Font.Color = $colorAmber
$skillCategory.Font.Size = 10

# --- SkillItems (paragraph) ---
$skillItems = $d.Styles.Add("SkillItems", 1)
$skillItems.NameLocal = "Skill Items"
$skillItems.ParagraphFormat.SpaceBefore = 0
$skillItems.ParagraphFormat.SpaceAfter = 0
$skillItems.ParagraphFormat.LineSpacingRule = 0
$skillItems.ParagraphFormat.LineSpacing = 14.4
$skillItems.Font.Name = "Liberation Serif"
$skillItems.Font.Bold = $false
$skillItems.Font.Color = $colorBrown
$skillItems.Font.Size = 10

# --- SkillHighlight (character) ---
$skillHighlight = $d.Styles.Add("SkillHighlight", 2)
$skillHighlight.NameLocal = "Skill Highlight"
$skillHighlight.Font.Name = "Liberation Serif"
$skillHighlight.Font.Bold = $true
$skillHighlight.Font.Color = $colorAmber
$skillHighlight.Font.Size = 10

# --- SkillLevel (character) ---
$skillLevel = $d.Styles.Add("SkillLevel", 2)
$skillLevel.NameLocal = "Skill Level"
$skillLevel.Font.Name = "Liberation Serif"
$skillLevel.Font.Bold = $false
$skillLevel.Font.Color = $colorAmber
$skillLevel.Font.Size = 10
